$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45188 to 45189
# for every data row (rows 2 through 556).
$ws.Range("C2:C556").Value = 45189
